$wb = $excel.ActiveWorkbook

# --- Shifts sheet ---
$shifts = $wb.Worksheets.Item("Shifts")
$shifts.Range("D5").Value = "SLC"

$shifts.Range("B6:D6").Copy()
$shifts.Range("B7").PasteSpecial(-4122)
$shifts.Range("B7").Value = $shifts.Range("B6").Value()
$shifts.Range("C7").Value = $shifts.Range("C6").Value()
$shifts.Range("D7").Value = "Reference"
$shifts.Range("D8").Select() | Out-Null

# --- Staff sheet ---
$staff = $wb.Worksheets.Item("Staff")
$staff.Range("X1").Value = "SLC?"
$staff.Range("Y1").Value = "Reference?"
$staff.Range("Z1").Value = "Standard? (empty means yes)"
$staff.Range("X1:Z1").Font.Bold = $true

$staff.Range("X2").Value = "Yes"
$staff.Range("Y2").Value = "y"
$staff.Range("Y3").Value = "No"
$staff.Range("Z3").Value = "n"

# --- Active sheet / selection ---
$staff.Activate()
$staff.Range("X3").Select() | Out-Null
